$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-15 03:06:47"
$wsZh.Range("G2").Value = "2016-01-15 03:07:27"

# "de-de" sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-15 03:06:58"
$wsDe.Range("G2").Value = "2016-01-15 03:07:42"
